$wb = $excel.ActiveWorkbook
$wsDatabase = $wb.Worksheets.Item("database")
$wsMetada = $wb.Worksheets.Item("metada")

# --- sheet1 (database): insert a new "wood" column before column X ---
$wsDatabase.Columns("X:X").Insert()
$wsDatabase.Range("X1").Value = "wood"
# match the column width that the inserted column should carry (same as column W)
$wsDatabase.Columns("X:X").ColumnWidth = $wsDatabase.Columns("W:W").ColumnWidth

# --- sheet2 (metada): insert a new "wood" / "amount of wood" row before row 24 ---
$wsMetada.Rows("24:24").Insert()
$wsMetada.Range("A24").Value = "wood"
$wsMetada.Range("B24").Value = "amount of wood"

# --- restore / update view state (selection, active cell, active sheet) ---
$wsDatabase.Activate()
$wsDatabase.Range("X1").Select()

$wsMetada.Activate()
$wsMetada.Range("F26").Select()
